$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is an updated event (values replaced); rows 3-5 are newly
# appended events from the daily export, extending the used range
# from A1:L2 to A1:L5.

# Row 2
$ws.Range("A2").Value = '281474991395157-1749154312868'
$ws.Range("B2").Value = 'Harsh Brake'
$ws.Range("C2").Value = '2025-06-05T14:11:52.868'
$ws.Range("D2").Value = "'281474991395157"
$ws.Range("E2").Value = "'126"
$ws.Range("F2").Value = 'No driver ID'
$ws.Range("G2").Value = 'No driver name'
$ws.Range("H2").Value = 20.65100792
$ws.Range("I2").Value = -103.3151144
$ws.Range("J2").Value = 0.7291494607925415
$ws.Range("K2").Value = 'No video URL'
$ws.Range("L2").Value = 'No video URL'

# Row 3
$ws.Range("A3").Value = '281474991205821-1749153333848'
$ws.Range("B3").Value = 'No Seat Belt'
$ws.Range("C3").Value = '2025-06-05T13:55:33.848'
$ws.Range("D3").Value = "'281474991205821"
$ws.Range("E3").Value = "'148"
$ws.Range("F3").Value = "'51834015"
$ws.Range("G3").Value = 'LUIS IBARRA'
$ws.Range("H3").Value = 20.643505709
$ws.Range("I3").Value = -103.35501921
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205821/1749153331348/OwQjboHWXD-camera-video-segment-driver-1749153333848.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOL7ZVC47%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T150627Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJIMEYCIQDTYgsp%2BIpTVyLfMQC9gcyy8bgjc7f%2F0gXRbRkORQg80QIhAJim9cI4G4RK23bUbuSlv91hpWpTlCsaNxB7mHk%2Bf3QGKt0DCF4QBBoMNzgxMjA0OTQyMjQ0IgwIarr7z0%2Fqn6KedwwqugNBPEun6%2FnAmsG1ghQpSyXrZ%2FrTgzqq0KX79U5D%2F62ZAbcibgxxUe%2Bc4QgDAwffuS6YSX8uWYcv4VIOTlT%2Fvwqlpd3LnsQXjvaproiOe6i7bviI6KMT4QVkdwIHDBWCgsyGsvWb1YETGg3T5eUwO1sxvoVmzx8kC7D2m41iW3Ki%2FfuRbIeH4lJ1N%2FfgHMvSf5hSxzgtSmcQ7laWLl0aocOEd%2FnB4soY6CxdZ0tqJTqtEIkke7ZFKqXDmgLQwPva2VMMI7TmQLJvVJCXAX0Zx6PW2Hq0IZ2tXSv5ogHkfQ%2FXpfIAAdZ9uSERgcOcQ%2B%2FgZ%2FIt%2FBNPWJ9kh5%2FfX4WgLP68UzbZFlWS6gGSp9diBjxCeRrygP9VAIINaOS2dmYG8kX%2F4FAyiZ5XDBjOLEEl0shTDIM6AtldOfTbvHvrP8gCzbWBo7GC5VayYQAHQNDQJ0ICMJQQzn2xc2ZrbRIkdCQmoapVZbKMS8LLdpiHKIwkWoBlDx7zBiQqMl7t5sV5SWElg9LCJ8srX3Nf8PlcjNpu1cFQM713sRcy2TunIqMy28dC55WWHqSRlOTFO2nABODH%2BsHd6tiHfXkUMJ7Ti8IGOqQBgtuHfPiQFFE%2F8N%2BCpwn9yEgo1qU1HsjSZAVIcUNThm0Zf60CY9G8OU9Ny0rvrJwzSpjSOYwyqEcTPNUmfIJ8ufTT%2BYzyFh59KnNd4tr3SNwrw8F%2BjGWRbfarGwgUdFzYXYJ4h8MNtjKu%2F3JEQFxslmG7yEsX3IsEAPSdPkMFFi2ppSlWdEMUdVc02I3jGD8n36KtJ%2FOQIBMLYmsVvYk5UnLK3qc%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2006%20Jun%202025%2023%3A06%3A27%20GMT&X-Amz-Signature=4d51d7d5bb136b2377bf83af1966a07c8945988f4f52f6eb940508868f4db2a6'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205821/1749153331348/oV5nVwJUOa-camera-video-segment-1749153333848.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOL7ZVC47%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T150627Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJIMEYCIQDTYgsp%2BIpTVyLfMQC9gcyy8bgjc7f%2F0gXRbRkORQg80QIhAJim9cI4G4RK23bUbuSlv91hpWpTlCsaNxB7mHk%2Bf3QGKt0DCF4QBBoMNzgxMjA0OTQyMjQ0IgwIarr7z0%2Fqn6KedwwqugNBPEun6%2FnAmsG1ghQpSyXrZ%2FrTgzqq0KX79U5D%2F62ZAbcibgxxUe%2Bc4QgDAwffuS6YSX8uWYcv4VIOTlT%2Fvwqlpd3LnsQXjvaproiOe6i7bviI6KMT4QVkdwIHDBWCgsyGsvWb1YETGg3T5eUwO1sxvoVmzx8kC7D2m41iW3Ki%2FfuRbIeH4lJ1N%2FfgHMvSf5hSxzgtSmcQ7laWLl0aocOEd%2FnB4soY6CxdZ0tqJTqtEIkke7ZFKqXDmgLQwPva2VMMI7TmQLJvVJCXAX0Zx6PW2Hq0IZ2tXSv5ogHkfQ%2FXpfIAAdZ9uSERgcOcQ%2B%2FgZ%2FIt%2FBNPWJ9kh5%2FfX4WgLP68UzbZFlWS6gGSp9diBjxCeRrygP9VAIINaOS2dmYG8kX%2F4FAyiZ5XDBjOLEEl0shTDIM6AtldOfTbvHvrP8gCzbWBo7GC5VayYQAHQNDQJ0ICMJQQzn2xc2ZrbRIkdCQmoapVZbKMS8LLdpiHKIwkWoBlDx7zBiQqMl7t5sV5SWElg9LCJ8srX3Nf8PlcjNpu1cFQM713sRcy2TunIqMy28dC55WWHqSRlOTFO2nABODH%2BsHd6tiHfXkUMJ7Ti8IGOqQBgtuHfPiQFFE%2F8N%2BCpwn9yEgo1qU1HsjSZAVIcUNThm0Zf60CY9G8OU9Ny0rvrJwzSpjSOYwyqEcTPNUmfIJ8ufTT%2BYzyFh59KnNd4tr3SNwrw8F%2BjGWRbfarGwgUdFzYXYJ4h8MNtjKu%2F3JEQFxslmG7yEsX3IsEAPSdPkMFFi2ppSlWdEMUdVc02I3jGD8n36KtJ%2FOQIBMLYmsVvYk5UnLK3qc%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2006%20Jun%202025%2023%3A06%3A27%20GMT&X-Amz-Signature=066405441171de6423f9b50312f0cdd6da2bc0f583e98459d1fc518634d17dc2'

# Row 4
$ws.Range("A4").Value = '281474991109254-1749147389968'
$ws.Range("B4").Value = 'No Seat Belt'
$ws.Range("C4").Value = '2025-06-05T12:16:29.968'
$ws.Range("D4").Value = "'281474991109254"
$ws.Range("E4").Value = "'117"
$ws.Range("F4").Value = "'51834027"
$ws.Range("G4").Value = 'JORGE HUMBERTO TORNERO'
$ws.Range("H4").Value = 20.68140625
$ws.Range("I4").Value = -103.46048785
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109254/1749147387468/GKgOu1ZeR3-camera-video-segment-driver-1749147389968.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOL7ZVC47%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T150627Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJIMEYCIQDTYgsp%2BIpTVyLfMQC9gcyy8bgjc7f%2F0gXRbRkORQg80QIhAJim9cI4G4RK23bUbuSlv91hpWpTlCsaNxB7mHk%2Bf3QGKt0DCF4QBBoMNzgxMjA0OTQyMjQ0IgwIarr7z0%2Fqn6KedwwqugNBPEun6%2FnAmsG1ghQpSyXrZ%2FrTgzqq0KX79U5D%2F62ZAbcibgxxUe%2Bc4QgDAwffuS6YSX8uWYcv4VIOTlT%2Fvwqlpd3LnsQXjvaproiOe6i7bviI6KMT4QVkdwIHDBWCgsyGsvWb1YETGg3T5eUwO1sxvoVmzx8kC7D2m41iW3Ki%2FfuRbIeH4lJ1N%2FfgHMvSf5hSxzgtSmcQ7laWLl0aocOEd%2FnB4soY6CxdZ0tqJTqtEIkke7ZFKqXDmgLQwPva2VMMI7TmQLJvVJCXAX0Zx6PW2Hq0IZ2tXSv5ogHkfQ%2FXpfIAAdZ9uSERgcOcQ%2B%2FgZ%2FIt%2FBNPWJ9kh5%2FfX4WgLP68UzbZFlWS6gGSp9diBjxCeRrygP9VAIINaOS2dmYG8kX%2F4FAyiZ5XDBjOLEEl0shTDIM6AtldOfTbvHvrP8gCzbWBo7GC5VayYQAHQNDQJ0ICMJQQzn2xc2ZrbRIkdCQmoapVZbKMS8LLdpiHKIwkWoBlDx7zBiQqMl7t5sV5SWElg9LCJ8srX3Nf8PlcjNpu1cFQM713sRcy2TunIqMy28dC55WWHqSRlOTFO2nABODH%2BsHd6tiHfXkUMJ7Ti8IGOqQBgtuHfPiQFFE%2F8N%2BCpwn9yEgo1qU1HsjSZAVIcUNThm0Zf60CY9G8OU9Ny0rvrJwzSpjSOYwyqEcTPNUmfIJ8ufTT%2BYzyFh59KnNd4tr3SNwrw8F%2BjGWRbfarGwgUdFzYXYJ4h8MNtjKu%2F3JEQFxslmG7yEsX3IsEAPSdPkMFFi2ppSlWdEMUdVc02I3jGD8n36KtJ%2FOQIBMLYmsVvYk5UnLK3qc%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2006%20Jun%202025%2023%3A06%3A27%20GMT&X-Amz-Signature=c37d388c8a19f6c8e16b79b7e3b0dbfe3820f170158d9145ae036e25107e9f7f'
$ws.Range("L4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109254/1749147387468/dQoeD1EUOr-camera-video-segment-1749147389968.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOL7ZVC47%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T150627Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIX%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJIMEYCIQDTYgsp%2BIpTVyLfMQC9gcyy8bgjc7f%2F0gXRbRkORQg80QIhAJim9cI4G4RK23bUbuSlv91hpWpTlCsaNxB7mHk%2Bf3QGKt0DCF4QBBoMNzgxMjA0OTQyMjQ0IgwIarr7z0%2Fqn6KedwwqugNBPEun6%2FnAmsG1ghQpSyXrZ%2FrTgzqq0KX79U5D%2F62ZAbcibgxxUe%2Bc4QgDAwffuS6YSX8uWYcv4VIOTlT%2Fvwqlpd3LnsQXjvaproiOe6i7bviI6KMT4QVkdwIHDBWCgsyGsvWb1YETGg3T5eUwO1sxvoVmzx8kC7D2m41iW3Ki%2FfuRbIeH4lJ1N%2FfgHMvSf5hSxzgtSmcQ7laWLl0aocOEd%2FnB4soY6CxdZ0tqJTqtEIkke7ZFKqXDmgLQwPva2VMMI7TmQLJvVJCXAX0Zx6PW2Hq0IZ2tXSv5ogHkfQ%2FXpfIAAdZ9uSERgcOcQ%2B%2FgZ%2FIt%2FBNPWJ9kh5%2FfX4WgLP68UzbZFlWS6gGSp9diBjxCeRrygP9VAIINaOS2dmYG8kX%2F4FAyiZ5XDBjOLEEl0shTDIM6AtldOfTbvHvrP8gCzbWBo7GC5VayYQAHQNDQJ0ICMJQQzn2xc2ZrbRIkdCQmoapVZbKMS8LLdpiHKIwkWoBlDx7zBiQqMl7t5sV5SWElg9LCJ8srX3Nf8PlcjNpu1cFQM713sRcy2TunIqMy28dC55WWHqSRlOTFO2nABODH%2BsHd6tiHfXkUMJ7Ti8IGOqQBgtuHfPiQFFE%2F8N%2BCpwn9yEgo1qU1HsjSZAVIcUNThm0Zf60CY9G8OU9Ny0rvrJwzSpjSOYwyqEcTPNUmfIJ8ufTT%2BYzyFh59KnNd4tr3SNwrw8F%2BjGWRbfarGwgUdFzYXYJ4h8MNtjKu%2F3JEQFxslmG7yEsX3IsEAPSdPkMFFi2ppSlWdEMUdVc02I3jGD8n36KtJ%2FOQIBMLYmsVvYk5UnLK3qc%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2006%20Jun%202025%2023%3A06%3A27%20GMT&X-Amz-Signature=15072135254e4e2e5b563e7af2873d8959413399c6d09eb7739a0f242c4eb94a'

# Row 5
$ws.Range("A5").Value = '281474991395097-1749136405801'
$ws.Range("B5").Value = 'Harsh Brake'
$ws.Range("C5").Value = '2025-06-05T09:13:25.801'
$ws.Range("D5").Value = "'281474991395097"
$ws.Range("E5").Value = "'125"
$ws.Range("F5").Value = "'51834055"
$ws.Range("G5").Value = 'DAVID SERRANO'
$ws.Range("H5").Value = 20.65729221
$ws.Range("I5").Value = -103.41939504
$ws.Range("J5").Value = 0.8426278233528137
$ws.Range("K5").Value = 'No video URL'
$ws.Range("L5").Value = 'No video URL'
